# Insert a new Title-styled paragraph " 020 simple chapter" before the
# existing first paragraph ("Chapter One", styled Heading1).
$d = $word.ActiveDocument

$firstPara = $d.Paragraphs(1)
$firstPara.Range.InsertParagraphBefore()

$titleRange = $d.Paragraphs(1).Range
$titleRange.Text = " 020 simple chapter"
$titleRange.Style = "Title"
